$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet1 (inventario): update existing Entradas/Salidas totals on row 2 ---
$ws1.Range("E2").Value = 24
$ws1.Range("F2").Value = 36

# --- Sheet1 (inventario): append new inventory rows 4 and 5 ---
$ws1.Range("A4").Value = "'power300"
$ws1.Range("A4").Style = "Normal"
$ws1.Range("B4").Value = "'antena"
$ws1.Range("B4").Style = "Normal"
$ws1.Range("C4").Value = 15
$ws1.Range("C4").Style = "Normal"
$ws1.Range("D4").Value = 0
$ws1.Range("D4").Style = "Normal"
$ws1.Range("E4").Value = 0
$ws1.Range("E4").Style = "Normal"
$ws1.Range("F4").Value = 15
$ws1.Range("F4").Style = "Normal"

$ws1.Range("A5").Value = "'185"
$ws1.Range("A5").Style = "Normal"
$ws1.Range("B5").Value = "'hola crayola"
$ws1.Range("B5").Style = "Normal"
$ws1.Range("C5").Value = 20
$ws1.Range("C5").Style = "Normal"
$ws1.Range("D5").Value = 0
$ws1.Range("D5").Style = "Normal"
$ws1.Range("E5").Value = 0
$ws1.Range("E5").Style = "Normal"
$ws1.Range("F5").Value = 20
$ws1.Range("F5").Style = "Normal"

# --- Sheet3 (Salidas): append new salida rows 4 through 10 ---
$ws3.Range("A4").Value = "'1852"
$ws3.Range("A4").Style = "Normal"
$ws3.Range("B4").Value = "'2018-11-09"
$ws3.Range("B4").Style = "Normal"
$ws3.Range("C4").Value = "'8965"
$ws3.Range("C4").Style = "Normal"
$ws3.Range("D4").Value = "'hola"
$ws3.Range("D4").Style = "Normal"
$ws3.Range("E4").Value = 4
$ws3.Range("E4").Style = "Normal"
$ws3.Range("F4").Value = "'sadas"
$ws3.Range("F4").Style = "Normal"
$ws3.Range("G4").Value = "'sadas"
$ws3.Range("G4").Style = "Normal"
$ws3.Range("H4").Value = "'sdadas"
$ws3.Range("H4").Style = "Normal"
$ws3.Range("I4").Value = "'6"
$ws3.Range("I4").Style = "Normal"

$ws3.Range("A5").Value = "'asd"
$ws3.Range("A5").Style = "Normal"
$ws3.Range("B5").Value = "'2018-11-17"
$ws3.Range("B5").Style = "Normal"
$ws3.Range("C5").Value = "'8965"
$ws3.Range("C5").Style = "Normal"
$ws3.Range("D5").Value = "'hola"
$ws3.Range("D5").Style = "Normal"
$ws3.Range("E5").Value = 5
$ws3.Range("E5").Style = "Normal"
$ws3.Range("F5").Value = "'sadasd"
$ws3.Range("F5").Style = "Normal"
$ws3.Range("G5").Value = "'400"
$ws3.Range("G5").Style = "Normal"
$ws3.Range("H5").Value = "'asdas"
$ws3.Range("H5").Style = "Normal"
$ws3.Range("I5").Value = "'5"
$ws3.Range("I5").Style = "Normal"

$ws3.Range("A6").Value = "'haber haber"
$ws3.Range("A6").Style = "Normal"
$ws3.Range("B6").Value = "'2018-11-10"
$ws3.Range("B6").Style = "Normal"
$ws3.Range("C6").Value = "'8965"
$ws3.Range("C6").Style = "Normal"
$ws3.Range("D6").Value = "'hola"
$ws3.Range("D6").Style = "Normal"
$ws3.Range("E6").Value = 3
$ws3.Range("E6").Style = "Normal"
$ws3.Range("F6").Value = "'lll"
$ws3.Range("F6").Style = "Normal"
$ws3.Range("G6").Value = "'333"
$ws3.Range("G6").Style = "Normal"
$ws3.Range("H6").Value = "'15"
$ws3.Range("H6").Style = "Normal"
$ws3.Range("I6").Value = "'4"
$ws3.Range("I6").Style = "Normal"

$ws3.Range("A7").Value = "'1852"
$ws3.Range("A7").Style = "Normal"
$ws3.Range("B7").Value = "'2018-11-03"
$ws3.Range("B7").Style = "Normal"
$ws3.Range("C7").Value = "'8965"
$ws3.Range("C7").Style = "Normal"
$ws3.Range("D7").Value = "'hola"
$ws3.Range("D7").Style = "Normal"
$ws3.Range("E7").Value = 3
$ws3.Range("E7").Style = "Normal"
$ws3.Range("F7").Value = "'Esta bien bonito"
$ws3.Range("F7").Style = "Normal"
$ws3.Range("G7").Value = "'sadasd"
$ws3.Range("G7").Style = "Normal"
$ws3.Range("H7").Value = "'sadas"
$ws3.Range("H7").Style = "Normal"
$ws3.Range("I7").Value = "'4"
$ws3.Range("I7").Style = "Normal"

$ws3.Range("A8").Value = "'asdas"
$ws3.Range("A8").Style = "Normal"
$ws3.Range("B8").Value = "'2018-11-16"
$ws3.Range("B8").Style = "Normal"
$ws3.Range("C8").Value = "'8965"
$ws3.Range("C8").Style = "Normal"
$ws3.Range("D8").Value = "'hola"
$ws3.Range("D8").Style = "Normal"
$ws3.Range("E8").Value = 3
$ws3.Range("E8").Style = "Normal"
$ws3.Range("F8").Value = "'sadasd"
$ws3.Range("F8").Style = "Normal"
$ws3.Range("G8").Value = "'sadas"
$ws3.Range("G8").Style = "Normal"
$ws3.Range("H8").Value = "'asdas"
$ws3.Range("H8").Style = "Normal"
$ws3.Range("I8").Value = "'3"
$ws3.Range("I8").Style = "Normal"

$ws3.Range("A9").Value = "'0123"
$ws3.Range("A9").Style = "Normal"
$ws3.Range("B9").Value = "'2018-11-02"
$ws3.Range("B9").Style = "Normal"
$ws3.Range("C9").Value = "'8965"
$ws3.Range("C9").Style = "Normal"
$ws3.Range("D9").Value = "'hola"
$ws3.Range("D9").Style = "Normal"
$ws3.Range("E9").Value = 1
$ws3.Range("E9").Style = "Normal"
$ws3.Range("F9").Value = "'instaLACION"
$ws3.Range("F9").Style = "Normal"
$ws3.Range("G9").Value = "'5253"
$ws3.Range("G9").Style = "Normal"
$ws3.Range("H9").Value = "'OSCAR"
$ws3.Range("H9").Style = "Normal"
$ws3.Range("I9").Value = "'01486513847"
$ws3.Range("I9").Style = "Normal"

$ws3.Range("A10").Value = "'0123"
$ws3.Range("A10").Style = "Normal"
$ws3.Range("B10").Value = "'2018-11-02"
$ws3.Range("B10").Style = "Normal"
$ws3.Range("C10").Value = "'8965"
$ws3.Range("C10").Style = "Normal"
$ws3.Range("D10").Value = "'hola"
$ws3.Range("D10").Style = "Normal"
$ws3.Range("E10").Value = 1
$ws3.Range("E10").Style = "Normal"
$ws3.Range("F10").Value = "'instaLACION"
$ws3.Range("F10").Style = "Normal"
$ws3.Range("G10").Value = "'ewg"
$ws3.Range("G10").Style = "Normal"
$ws3.Range("H10").Value = "'OSCAR"
$ws3.Range("H10").Style = "Normal"
$ws3.Range("I10").Value = "'23131614510230"
$ws3.Range("I10").Style = "Normal"
